# Rename the four stack columns with descriptive subsystem names and
# retune the D/E stack-size formulas (free stacks & heap monitoring).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Stack2 - SDCard"
$ws.Range("D1").Value = "Stack1 - SysMon"
$ws.Range("F1").Value = "Stack3 - WlsCom"
$ws.Range("G1").Value = "Stack4 - IrrCtrl"

# Widen the newly-relevant columns so the longer labels aren't truncated.
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 14.666666666666666
$ws.Columns.Item(6).ColumnWidth = 15.833333333333334
$ws.Columns.Item(7).ColumnWidth = 13.0

# Corrected multipliers used to estimate stack usage (publishLog fix).
$ws.Range("D2").Formula = "=30*B2+C2"
$ws.Range("E2").Formula = "=80*B2+C2"

$ws.Range("H6").Select()
